$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-23 04:47:48"
$wsZh.Range("H2").Value = "2016-03-23 04:48:30"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-23 04:47:57"
$wsDe.Range("H2").Value = "2016-03-23 04:48:46"
